# FamilyComposition.docx refactor:
#  1. Merge the "Видана ... Сумської області " / "про те, що " runs into one run.
#  2. Turn " зареєстрованій за адресою: ${person_address}" into
#     "  ${person_address_registration}" split as three runs:
#     "  ${person_address" / "_registration" / "}".
#  3. Give the empty run in the blank paragraph right after "${/relatives}"
#     the same run formatting (bCs / sz28 / szCs28) as its neighbours.

$d = $word.ActiveDocument

# --- 1) merge "Видана ... область " + "про те, що " into a single run ---------
$merge = $d.Content
$mergeText = "Видана   виконавчим комітетом Попівської сільської ради Конотопського району Сумської області про те, що "
$merge.Find.ClearFormatting()
[void]$merge.Find.Execute($mergeText, $false, $false, $false, $false, $false, $true, 1, $false, $mergeText, 2)

# --- 2) rewrite the registration-address sentence fragment --------------------
$old2 = " зареєстрованій за адресою: `$`{person_address`}"
$new2 = "  `$`{person_address"
$r2 = $d.Content
$r2.Find.ClearFormatting()
[void]$r2.Find.Execute($old2, $false, $false, $false, $false, $false, $true, 1, $false, $new2, 2)

# locate the spot right after "...${person_address" that Find just produced
$r2b = $d.Content
$r2b.Find.ClearFormatting()
[void]$r2b.Find.Execute($new2, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertPos = $r2b.End
$ins = $d.Range($insertPos, $insertPos)
$ins.InsertAfter("_registration}")

# force the freshly inserted text to live in its own run(s) even though its
# formatting matches the surrounding text, by toggling a character property
# on then off again (mirrors what Word does when a run boundary is created
# by an intermediate formatting change).
$splitLen = ("_registration}").Length
$splitRange = $d.Range($insertPos, $insertPos + $splitLen)
$splitRange.Font.Bold = $true
$splitRange.Font.Bold = $false

# split "_registration" away from the trailing "}" the same way
$regLen = ("_registration").Length
$regRange = $d.Range($insertPos, $insertPos + $regLen)
$regRange.Font.Bold = $true
$regRange.Font.Bold = $false

# --- 3) add run formatting to the empty run after "${/relatives}" -------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "`$`{/relatives`}`r") {
        $blank = $d.Paragraphs.Item($i + 1)
        $blankRange = $blank.Range
        $blankRange.Font.BoldBi = $true
        $blankRange.Font.Size = 14
        $blankRange.Font.SizeBi = 14
        break
    }
}

Write-Output "done"
